{"js": "const table = context.document.body.tables.getFirst();\n\nconst replacements = [\n  { row: 0, col: 0, text: \"91\u00f78=\" },\n  { row: 0, col: 1, text: \"95\u00f72=\" },\n  { row: 0, col: 2, text: \"15\u00f76=\" },\n  { row: 0, col: 3, text: \"50\u00f79=\" },\n  { row: 0, col: 4, text: \"88\u00f76=\" },\n  { row: 4, col: 0, text: \"16\u00f73=\" },\n  { row: 4, col: 1, text: \"40\u00f76=\" },\n  { row: 4, col: 2, text: \"42\u00f74=\" },\n  { row: 4, col: 3, text: \"35\u00f76=\" },\n  { row: 4, col: 4, text: \"67\u00f72=\" },\n  { row: 8, col: 0, text: \"56\u00f78=\" },\n  { row: 8, col: 1, text: \"26\u00f73=\" },\n  { row: 8, col: 2, text: \"23\u00f75=\" },\n  { row: 8, col: 3, text: \"26\u00f77=\" },\n  { row: 8, col: 4, text: \"12\u00f75=\" },\n  { row: 12, col: 0, text: \"42\u00f76=\" },\n  { row: 12, col: 1, text: \"28\u00f78=\" },\n  { row: 12, col: 2, text: \"49\u00f76=\" },\n  { row: 12, col: 3, text: \"34\u00f72=\" },\n  { row: 12, col: 4, text: \"48\u00f75=\" },\n  { row: 16, col: 0, text: \"24\u00f75=\" },\n  { row: 16, col: 1, text: \"86\u00f79=\" },\n  { row: 16, col: 2, text: \"21\u00f77=\" },\n  { row: 16, col: 3, text: \"94\u00f73=\" },\n  { row: 16, col: 4, text: \"56\u00f79=\" },\n];\n\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  const para = cell.body.paragraphs.getFirst();\n  para.insertText(r.text, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; Text = \"91\u00f78=\" }\n    @{ Row = 1; Col = 2; Text = \"95\u00f72=\" }\n    @{ Row = 1; Col = 3; Text = \"15\u00f76=\" }\n    @{ Row = 1; Col = 4; Text = \"50\u00f79=\" }\n    @{ Row = 1; Col = 5; Text = \"88\u00f76=\" }\n    @{ Row = 5; Col = 1; Text = \"16\u00f73=\" }\n    @{ Row = 5; Col = 2; Text = \"40\u00f76=\" }\n    @{ Row = 5; Col = 3; Text = \"42\u00f74=\" }\n    @{ Row = 5; Col = 4; Text = \"35\u00f76=\" }\n    @{ Row = 5; Col = 5; Text = \"67\u00f72=\" }\n    @{ Row = 9; Col = 1; Text = \"56\u00f78=\" }\n    @{ Row = 9; Col = 2; Text = \"26\u00f73=\" }\n    @{ Row = 9; Col = 3; Text = \"23\u00f75=\" }\n    @{ Row = 9; Col = 4; Text = \"26\u00f77=\" }\n    @{ Row = 9; Col = 5; Text = \"12\u00f75=\" }\n    @{ Row = 13; Col = 1; Text = \"42\u00f76=\" }\n    @{ Row = 13; Col = 2; Text = \"28\u00f78=\" }\n    @{ Row = 13; Col = 3; Text = \"49\u00f76=\" }\n    @{ Row = 13; Col = 4; Text = \"34\u00f72=\" }\n    @{ Row = 13; Col = 5; Text = \"48\u00f75=\" }\n    @{ Row = 17; Col = 1; Text = \"24\u00f75=\" }\n    @{ Row = 17; Col = 2; Text = \"86\u00f79=\" }\n    @{ Row = 17; Col = 3; Text = \"21\u00f77=\" }\n    @{ Row = 17; Col = 4; Text = \"94\u00f73=\" }\n    @{ Row = 17; Col = 5; Text = \"56\u00f79=\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $cell.Range.Text = $r.Text\n}"}
